# "forgot to turn off vpn redid renton data"
# Re-measured the "Renton, WA" latency samples (column I, rows 2-21) after
# realizing the earlier run was taken over a VPN. The AVERAGE formula in I22
# recalculates automatically. Also nudges the active selection, matching
# where the user's cursor ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rentonLatency = @{
    2  = 11.51
    3  = 14.67
    4  = 13.38
    5  = 14.08
    6  = 13.78
    7  = 14.32
    8  = 13.46
    9  = 13.31
    10 = 13.83
    11 = 13.31
    12 = 13.83
    13 = 13.07
    14 = 13.95
    15 = 13.38
    16 = 12.86
    17 = 13.26
    18 = 13.3
    19 = 13.96
    20 = 14
    21 = 13.8
}

foreach ($row in $rentonLatency.Keys) {
    $ws.Cells.Item($row, 9).Value = $rentonLatency[$row]
}

$ws.Range("D18").Select()
